# Insert a new weekly record as row 5 (pushing the existing rows 5-30 down to 6-31)
# and populate it with the new "Provincia de Limarí" observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("5:5").Insert()

$ws.Cells.Item(5, 1).Value  = 11
$ws.Cells.Item(5, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(5, 3).Value  = "Bíobío"
$ws.Cells.Item(5, 4).Value  = 44532
$ws.Cells.Item(5, 5).Value  = 8
$ws.Cells.Item(5, 6).Value  = 100112031
$ws.Cells.Item(5, 7).Value  = "Poroto verde"
$ws.Cells.Item(5, 8).Value  = "Magnum"
$ws.Cells.Item(5, 9).Value  = "Primera"
$ws.Cells.Item(5, 10).Value = 250
$ws.Cells.Item(5, 11).Value = 33000
$ws.Cells.Item(5, 12).Value = 35000
$ws.Cells.Item(5, 13).Value = 33800
$ws.Cells.Item(5, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(5, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(5, 16).Value = 1352
$ws.Cells.Item(5, 17).Value = 25
$ws.Cells.Item(5, 18).Value = "Hortaliza"
